$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 56.216147
$ws.Range("H2").Value = 168.648441
$ws.Range("I2").Value = 0.9695233148109074
$ws.Range("J2").Value = 0.9695233148109074
$ws.Range("M2").Value = 1.168007333333333
$ws.Range("N2").Value = 3.504022
$ws.Range("O2").Value = 0.1638609704511517
$ws.Range("P2").Value = 0.1638609704511517
$ws.Range("Q2").Value = 65.66087194774467
$ws.Range("R2").Value = 590.947847529702
$ws.Range("S2").Value = 0.1588670312399328
$ws.Range("T2").Value = 0.1588670312399328
$ws.Range("G3").Value = 56.216147
$ws.Range("H3").Value = 168.648441
$ws.Range("I3").Value = 0.9695233148109074
$ws.Range("J3").Value = 0.9695233148109074
$ws.Range("O3").Value = 0.5019752511630595
$ws.Range("P3").Value = 0.5019752511630595
$ws.Range("Q3").Value = 201.146939364553
$ws.Range("R3").Value = 1810.322454280977
$ws.Range("S3").Value = 0.4866767094606472
$ws.Range("T3").Value = 0.4866767094606472
$ws.Range("G4").Value = 56.216147
$ws.Range("H4").Value = 168.648441
$ws.Range("I4").Value = 0.9695233148109074
$ws.Range("J4").Value = 0.9695233148109074
$ws.Range("M4").Value = 1.915392333333333
$ws.Range("N4").Value = 5.746177
$ws.Range("O4").Value = 0.2687123938160456
$ws.Range("P4").Value = 0.2687123938160456
$ws.Range("Q4").Value = 107.6759769733397
$ws.Range("R4").Value = 969.083792760057
$ws.Range("S4").Value = 0.2605229307833065
$ws.Range("T4").Value = 0.2605229307833065
$ws.Range("G5").Value = 56.216147
$ws.Range("H5").Value = 168.648441
$ws.Range("I5").Value = 0.9695233148109074
$ws.Range("J5").Value = 0.9695233148109074
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.46654
$ws.Range("N5").Value = 1.39962
$ws.Range("O5").Value = 0.06545138456974327
$ws.Range("P5").Value = 0.06545138456974327
$ws.Range("Q5").Value = 26.22708122138
$ws.Range("R5").Value = 236.04373099242
$ws.Range("S5").Value = 0.06345664332702097
$ws.Range("T5").Value = 0.06345664332702097
$ws.Range("I6").Value = 0.001287520467967504
$ws.Range("J6").Value = 0.001287520467967504
$ws.Range("M6").Value = 1.168007333333333
$ws.Range("N6").Value = 3.504022
$ws.Range("O6").Value = 0.1638609704511517
$ws.Range("P6").Value = 0.1638609704511517
$ws.Range("Q6").Value = 0.08719719813422222
$ws.Range("R6").Value = 0.784774783208
$ws.Range("S6").Value = 0.0002109743533568763
$ws.Range("T6").Value = 0.0002109743533568763
$ws.Range("I7").Value = 0.001287520467967504
$ws.Range("J7").Value = 0.001287520467967504
$ws.Range("O7").Value = 0.5019752511630595
$ws.Range("P7").Value = 0.5019752511630595
$ws.Range("S7").Value = 0.0006463034102855679
$ws.Range("T7").Value = 0.000646303410285568
$ws.Range("I8").Value = 0.001287520467967504
$ws.Range("J8").Value = 0.001287520467967504
$ws.Range("M8").Value = 1.915392333333333
$ws.Range("N8").Value = 5.746177
$ws.Range("O8").Value = 0.2687123938160456
$ws.Range("P8").Value = 0.2687123938160456
$ws.Range("Q8").Value = 0.1429929761808889
$ws.Range("R8").Value = 1.286936785628
$ws.Range("S8").Value = 0.0003459727070347033
$ws.Range("T8").Value = 0.0003459727070347034
$ws.Range("I9").Value = 0.001287520467967504
$ws.Range("J9").Value = 0.001287520467967504
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.46654
$ws.Range("N9").Value = 1.39962
$ws.Range("O9").Value = 0.06545138456974327
$ws.Range("P9").Value = 0.06545138456974327
$ws.Range("Q9").Value = 0.03482938818666666
$ws.Range("R9").Value = 0.31346449368
$ws.Range("S9").Value = 0.00008426999729035693
$ws.Range("T9").Value = 0.00008426999729035695
$ws.Range("G10").Value = 0.4660483333333333
$ws.Range("H10").Value = 1.398145
$ws.Range("I10").Value = 0.008037632408272877
$ws.Range("J10").Value = 0.008037632408272877
$ws.Range("M10").Value = 1.168007333333333
$ws.Range("N10").Value = 3.504022
$ws.Range("O10").Value = 0.1638609704511517
$ws.Range("P10").Value = 0.1638609704511517
$ws.Range("Q10").Value = 0.5443478710211112
$ws.Range("R10").Value = 4.89913083919
$ws.Range("S10").Value = 0.001317054246549221
$ws.Range("T10").Value = 0.001317054246549221
$ws.Range("G11").Value = 0.4660483333333333
$ws.Range("H11").Value = 1.398145
$ws.Range("I11").Value = 0.008037632408272877
$ws.Range("J11").Value = 0.008037632408272877
$ws.Range("O11").Value = 0.5019752511630595
$ws.Range("P11").Value = 0.5019752511630595
$ws.Range("Q11").Value = 1.667567075451666
$ws.Range("R11").Value = 15.008103679065
$ws.Range("S11").Value = 0.004034692546899124
$ws.Range("T11").Value = 0.004034692546899124
$ws.Range("G12").Value = 0.4660483333333333
$ws.Range("H12").Value = 1.398145
$ws.Range("I12").Value = 0.008037632408272877
$ws.Range("J12").Value = 0.008037632408272877
$ws.Range("M12").Value = 1.915392333333333
$ws.Range("N12").Value = 5.746177
$ws.Range("O12").Value = 0.2687123938160456
$ws.Range("P12").Value = 0.2687123938160456
$ws.Range("Q12").Value = 0.8926654046294444
$ws.Range("R12").Value = 8.033988641665001
$ws.Range("S12").Value = 0.002159811445040432
$ws.Range("T12").Value = 0.002159811445040432
$ws.Range("G13").Value = 0.4660483333333333
$ws.Range("H13").Value = 1.398145
$ws.Range("I13").Value = 0.008037632408272877
$ws.Range("J13").Value = 0.008037632408272877
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.46654
$ws.Range("N13").Value = 1.39962
$ws.Range("O13").Value = 0.06545138456974327
$ws.Range("P13").Value = 0.06545138456974327
$ws.Range("Q13").Value = 0.2174301894333333
$ws.Range("R13").Value = 1.9568717049
$ws.Range("S13").Value = 0.0005260741697840998
$ws.Range("T13").Value = 0.0005260741697840998
$ws.Range("G14").Value = 1.226435333333333
$ws.Range("H14").Value = 3.679306
$ws.Range("I14").Value = 0.02115153231285227
$ws.Range("J14").Value = 0.02115153231285228
$ws.Range("M14").Value = 1.168007333333333
$ws.Range("N14").Value = 3.504022
$ws.Range("O14").Value = 0.1638609704511517
$ws.Range("P14").Value = 0.1638609704511517
$ws.Range("Q14").Value = 1.432485463192444
$ws.Range("R14").Value = 12.892369168732
$ws.Range("S14").Value = 0.003465910611312867
$ws.Range("T14").Value = 0.003465910611312868
$ws.Range("G15").Value = 1.226435333333333
$ws.Range("H15").Value = 3.679306
$ws.Range("I15").Value = 0.02115153231285227
$ws.Range("J15").Value = 0.02115153231285228
$ws.Range("O15").Value = 0.5019752511630595
$ws.Range("P15").Value = 0.5019752511630595
$ws.Range("Q15").Value = 4.388307039764666
$ws.Range("R15").Value = 39.49476335788199
$ws.Range("S15").Value = 0.01061754574522759
$ws.Range("T15").Value = 0.01061754574522759
$ws.Range("G16").Value = 1.226435333333333
$ws.Range("H16").Value = 3.679306
$ws.Range("I16").Value = 0.02115153231285227
$ws.Range("J16").Value = 0.02115153231285228
$ws.Range("M16").Value = 1.915392333333333
$ws.Range("N16").Value = 5.746177
$ws.Range("O16").Value = 0.2687123938160456
$ws.Range("P16").Value = 0.2687123938160456
$ws.Range("Q16").Value = 2.349104834795778
$ws.Range("R16").Value = 21.141943513162
$ws.Range("S16").Value = 0.005683678880663972
$ws.Range("T16").Value = 0.005683678880663975
$ws.Range("G17").Value = 1.226435333333333
$ws.Range("H17").Value = 3.679306
$ws.Range("I17").Value = 0.02115153231285227
$ws.Range("J17").Value = 0.02115153231285228
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.46654
$ws.Range("N17").Value = 1.39962
$ws.Range("O17").Value = 0.06545138456974327
$ws.Range("P17").Value = 0.06545138456974327
$ws.Range("Q17").Value = 0.5721811404133332
$ws.Range("R17").Value = 5.14963026372
$ws.Range("S17").Value = 0.001384397075647845
$ws.Range("T17").Value = 0.001384397075647846

Write-Output "Done applying changes"
